$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (from an existing header cell, e.g. A1)
# onto the three new header cells so they reuse the same style index
# instead of creating new (duplicate) style entries.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record values repeated for every data row (2-48)
$ws.Range("AD2:AD48").Value = 67
$ws.Range("AE2:AE48").Value = 95
$ws.Range("AF2:AF48").Value = 0
